$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values that changed
$ws.Range("O2").Value = 0.1810521476743106
$ws.Range("P2").Value = 0.1810521476743105
$ws.Range("S2").Value = 0.1810521476743106
$ws.Range("T2").Value = 0.1810521476743105

# Add new row 3 with data for FAPs -> MuSCs (Wnt1/Fzd10)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt1"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.422259
$ws.Range("H3").Value = 1.266777
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05581866666666666
$ws.Range("N3").Value = 0.167456
$ws.Range("O3").Value = 0.8189478523256895
$ws.Range("P3").Value = 0.8189478523256895
$ws.Range("Q3").Value = 0.023569934368
$ws.Range("R3").Value = 0.212129409312
$ws.Range("S3").Value = 0.8189478523256895
$ws.Range("T3").Value = 0.8189478523256895
